Write-Output "before"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(41)
$range = $s.Shapes.Range(@(3,4,5,6,7,8,9))
$grp = $range.Group()
Write-Output "grouped: $($grp.Name) id=$($grp.Id)"
Write-Output "Left=$($grp.Left) Top=$($grp.Top) Width=$($grp.Width) Height=$($grp.Height)"
$grp.Left = 252.0
$grp.Top = 138.08669291338583
$grp.Width = 420.0
$grp.Height = 271.95511811023624
Write-Output "after: Left=$($grp.Left) Top=$($grp.Top) Width=$($grp.Width) Height=$($grp.Height)"
